$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "326.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.71%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.05%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.250"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.97%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08374"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.26%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.936"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.85%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9722"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.06%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1153"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.81%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1901"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.95%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09708"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.41%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04625"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.16%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1058"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.41%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001296"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.90%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005812"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.48%"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.85%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.451"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.55%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3361"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.68%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.641"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-14.89%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.94%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2583"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.56%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04157"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.57%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001236"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.16%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004423"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.80%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001305"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.11%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002989"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-19.93%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02718"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "1.64%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05635"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.11%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007866"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.35%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1411"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.11%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007341"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.47%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002047"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.71%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007885"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.96%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3498"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006845"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.26%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.53%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003498"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.65%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003542"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "40.82%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.53%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.53%"
